$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-07 06:47:51"
$ws.Range("O2").Value = "-1.8 °C"
$ws.Range("E3").Value = "2026-02-07 06:47:53"
$ws.Range("E4").Value = "2026-02-07 06:47:56"
$ws.Range("H4").Value = "58%"
$ws.Range("J4").Value = "1001.4 hPa"
$ws.Range("N4").Value = "9.7 °C 6:24 TU"
$ws.Range("O4").Value = "11.2 °C"
$ws.Range("E5").Value = "2026-02-07 06:47:58"
$ws.Range("J5").Value = "1001.4 hPa"
$ws.Range("O5").Value = "8.4 °C"
$ws.Range("E6").Value = "2026-02-07 06:48:01"
$ws.Range("J6").Value = "1003.0 hPa"
$ws.Range("E7").Value = "2026-02-07 06:48:03"
$ws.Range("L7").Value = "45.0 km/h - 41º 6:09 TU"
$ws.Range("E8").Value = "2026-02-07 06:48:06"
$ws.Range("H8").Value = "94%"
$ws.Range("O8").Value = "4.0 °C"
$ws.Range("E9").Value = "2026-02-07 06:48:08"
$ws.Range("O9").Value = "1.4 °C"
$ws.Range("E10").Value = "2026-02-07 06:48:11"
$ws.Range("E11").Value = "2026-02-07 06:48:13"
$ws.Range("J11").Value = "1005.7 hPa"
$ws.Range("E12").Value = "2026-02-07 06:48:15"
$ws.Range("H12").Value = "70%"
$ws.Range("M12").Value = "11.3 °C 6:06 TU"
$ws.Range("O12").Value = "9.7 °C"
$ws.Range("E13").Value = "2026-02-07 06:48:18"
$ws.Range("H13").Value = "86%"
$ws.Range("O13").Value = "7.9 °C"
$ws.Range("E14").Value = "2026-02-07 06:48:20"
$ws.Range("I14").Value = "0.1 mm"
$ws.Range("N14").Value = "-7.8 °C 6:27 TU"
$ws.Range("O14").Value = "-5.8 °C"
$ws.Range("E15").Value = "2026-02-07 06:48:23"
$ws.Range("H15").Value = "88%"
$ws.Range("J15").Value = "1001.7 hPa"
$ws.Range("N15").Value = "2.6 °C 6:14 TU"
$ws.Range("O15").Value = "5.9 °C"
$ws.Range("E16").Value = "2026-02-07 06:48:25"
$ws.Range("N16").Value = "1.8 °C 6:14 TU"
$ws.Range("O16").Value = "2.8 °C"
$ws.Range("E17").Value = "2026-02-07 06:48:28"
$ws.Range("J17").Value = "1005.0 hPa"
$ws.Range("N17").Value = "2.3 °C 6:20 TU"
$ws.Range("O17").Value = "3.2 °C"
$ws.Range("E18").Value = "2026-02-07 06:48:31"
$ws.Range("O18").Value = "-7.7 °C"
$ws.Range("E19").Value = "2026-02-07 06:48:33"
$ws.Range("J19").Value = "1006.3 hPa"
$ws.Range("N19").Value = "1.5 °C 6:16 TU"
$ws.Range("O19").Value = "4.1 °C"
$ws.Range("E20").Value = "2026-02-07 06:48:36"
$ws.Range("H20").Value = "82%"
$ws.Range("E21").Value = "2026-02-07 06:48:38"
$ws.Range("H21").Value = "77%"
$ws.Range("J21").Value = "1002.1 hPa"
$ws.Range("O21").Value = "6.1 °C"
$ws.Range("E22").Value = "2026-02-07 06:48:41"
$ws.Range("H22").Value = "90%"
$ws.Range("L22").Value = "15.8 km/h - 358º 6:17 TU"
$ws.Range("O22").Value = "6.3 °C"
$ws.Range("E23").Value = "2026-02-07 06:48:43"
$ws.Range("J23").Value = "1001.5 hPa"
$ws.Range("N23").Value = "6.4 °C 6:29 TU"
$ws.Range("O23").Value = "7.5 °C"
$ws.Range("E24").Value = "2026-02-07 06:48:46"
$ws.Range("J24").Value = "1000.8 hPa"
$ws.Range("L24").Value = "68.8 km/h - 343º 6:10 TU"
$ws.Range("N24").Value = "9.7 °C 6:26 TU"
$ws.Range("E25").Value = "2026-02-07 06:48:49"
$ws.Range("H25").Value = "97%"
$ws.Range("J25").Value = "1005.3 hPa"
$ws.Range("E26").Value = "2026-02-07 06:48:51"
$ws.Range("O26").Value = "-2.3 °C"
$ws.Range("E27").Value = "2026-02-07 06:48:54"
$ws.Range("J27").Value = "1001.3 hPa"
$ws.Range("E28").Value = "2026-02-07 06:48:57"
$ws.Range("E29").Value = "2026-02-07 06:48:59"
$ws.Range("N29").Value = "8.6 °C 6:29 TU"
$ws.Range("O29").Value = "10.9 °C"
$ws.Range("E30").Value = "2026-02-07 06:49:02"
$ws.Range("E31").Value = "2026-02-07 06:49:04"
$ws.Range("J31").Value = "1005.9 hPa"
$ws.Range("E32").Value = "2026-02-07 06:49:06"
$ws.Range("H32").Value = "58%"
$ws.Range("J32").Value = "1004.5 hPa"
$ws.Range("E33").Value = "2026-02-07 06:49:09"
$ws.Range("E34").Value = "2026-02-07 06:49:12"
$ws.Range("H34").Value = "78%"
$ws.Range("K34").Value = "-0.1 MJ/m2"
$ws.Range("N34").Value = "4.7 °C 6:16 TU"
$ws.Range("O34").Value = "6.2 °C"
$ws.Range("E35").Value = "2026-02-07 06:49:14"
$ws.Range("O35").Value = "-6.0 °C"
$ws.Range("E36").Value = "2026-02-07 06:49:17"
$ws.Range("J36").Value = "1006.8 hPa"
$ws.Range("K36").Value = "-0.1 MJ/m2"
$ws.Range("N36").Value = "3.7 °C 6:29 TU"
